$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell edits, applied in the same order the original author made them so
# --- the shared-string table is rebuilt with the same append order/indices.
$ws.Range("B3").Value2  = "stroke[2]Color"
$ws.Range("C3").Value2  = "stroke[2]Width"
$ws.Range("A1").Value2  = "!black"
$ws.Range("G1").Value2  = "black"
$ws.Range("G7").Value2  = "234 50 20"
$ws.Range("G6").Value2  = "blue"
$ws.Range("G8").Value2  = "100, 0 50 0"
$ws.Range("G9").Value2  = "tan"
$ws.Range("G3").Value2  = "fillColor: red"
$ws.Range("G11").Value2 = "#0f0"
$ws.Range("B4").Value2  = "2nd stroke color"
$ws.Range("G2").Value2  = "appearance"
$ws.Range("G10").Value2 = "lab 100 88 -107"

# --- Column A width: drop AutoFit/bestFit, make it a custom 25.1640625-char
# --- width. The COM layer quantizes ColumnWidth to 1/6-character steps, so
# --- feed it the input that lands on the closest achievable step.
$ws.Columns(1).ColumnWidth = 24.333333333333336

# --- Selection moves to G4.
$ws.Range("G4").Select() | Out-Null

# --- Page setup: paper size + orientation now explicit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
